$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing columns (A->B, B->C, C->D, D->E)
$ws.Columns.Item(1).Insert()

# New column A header/value
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# Updated Neo4j queries (columns B2/C2, which used to be the old A2/B2 before the insert)
# Single-quoted here-strings (@'...'@) are used so backticks/$ are kept verbatim.
$caseQuery = @'
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "WHITE"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '') AS `Trial Code`,
    COALESCE(a.arm_id, '') AS `Arm`,
    COALESCE(a.arm_drug, '') AS `Arm Treatment`,
    COALESCE(c.disease, '') AS `Diagnosis`,
    COALESCE(c.gender, '') AS `Gender`,
    COALESCE(c.race, '') AS `Race`,
    COALESCE(c.ethnicity, '') AS `Ethnicity`
'@

$statQuery = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "WHITE"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
'@

$ws.Range("B2").Value = $caseQuery
$ws.Range("C2").Value = $statQuery

# Row 2 height grew to fit the longer wrapped query text
$ws.Rows.Item(2).RowHeight = 174

# Column widths: A is now a narrow "bestFit" index column; B/C/D/E keep the prior
# (unchanged) widths automatically after the column insert, so only A needs sizing.
$ws.Columns.Item(1).ColumnWidth = 8

# Selection moved
$ws.Range("B10").Select() | Out-Null
